# New adata function: display
# Rework of the "survey" sheet: the note rows that render calculated values
# now use richer HTML markup (bold / colored font), the stray if/else test
# scaffolding (rows 12-25 in the old layout) is removed, and the age-display
# note that used to live further down the sheet now directly follows the
# "display" note (new row 12), with the sheet ending at row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Wipe everything from row 12 down (the old if/else scaffolding and the
# trailing rows); we'll rewrite the small amount that remains afterwards.
$ws.Range("A12:M26").Clear()

# "display" note: drop the space after the colon.
$ws.Range("G9").Value = "display:{{calculates.display}}"

# Update the wording of the "diff in days" / "diff in years" note rows so
# the calculated values are wrapped in markup.
$ws.Range("G7").Value = "diff in days: <b>{{calculates.diffdays}}</b>"
$ws.Range("H7").Value = "diff in days: <b>{{calculates.diffdays}}</b>"

$ws.Range("G8").Value = "diff in years: <font color=""red"">{{calculates.diffInYears}}</font>"
$ws.Range("H8").Value = "diff in years: <font color=""red"">{{calculates.diffInYears}}</font>"

# The age-display note (previously further down the sheet, inside the old
# if/else scaffolding) now directly follows the "display" note.
$ws.Range("D12").Value = "note"
$ws.Range("G12").Value = "Alder I hele år: <b>{{calculates.ageInYears}}</b><br/>I måneder: {{calculates.ageInMonths}}<br/>I dage: {{calculates.ageInDays}}<br/>( baseret på {{data.ADA}} )<br/>"
$ws.Range("H12").Value = "Alder I hele år: <b>{{calculates.ageInYears}}</b><br/>I måneder: {{calculates.ageInMonths}}<br/>I dage: {{calculates.ageInDays}}<br/>( baseret på {{data.ADA}} )<br/>"

$ws.Range("B13").Value = "end screen"

# Tidy up the sheet view so it matches the new, much shorter, layout.
$ws.Range("G8").Select()
